# Capitalize the section-header labels (Russian/English) that were
# previously lower-case, and clear the lingering cell selection that
# pointed at A23 so the sheet opens with the default A1 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# "по полу / by sex" header row -> capitalised
$ws.Range("B14").Value = "По полу"
$ws.Range("C14").Value = "By sex"
$ws.Range("A14").Value = "Жынысы боюнча"

# "образование матери / education of mother" header row -> capitalised
$ws.Range("A17").Value = "Энесинин билими "
$ws.Range("B17").Value = "Образование матери "
$ws.Range("C17").Value = "Education of mother"

# "квинтиль по индексу благосостояния / wealth quintile" header row -> capitalised
$ws.Range("B23").Value = "Квинтиль по индексу благосостояния"
$ws.Range("C23").Value = "Wealth quintile"

# Reset the saved selection back to A1 (removes the stored A23 selection)
$ws.Range("A1").Select()
